{"js": "// Replace the 25 \"divisor problem\" answers in the document's table with\n// their newly generated values, cell by cell, in document order.\n// Each entry is [rowIndex, colIndex, oldText, newText] (0-based).\nconst replacements = [\n  [0, 0, \"809\u00f79=89, 8\", \"749\u00f76=124, 5\"],\n  [0, 1, \"271\u00f77=38, 5\", \"115\u00f79=12, 7\"],\n  [0, 2, \"841\u00f78=105, 1\", \"673\u00f74=168, 1\"],\n  [0, 3, \"303\u00f76=50, 3\", \"847\u00f79=94, 1\"],\n  [0, 4, \"886\u00f74=221, 2\", \"702\u00f76=117, 0\"],\n  [4, 0, \"845\u00f79=93, 8\", \"176\u00f78=22, 0\"],\n  [4, 1, \"324\u00f73=108, 0\", \"342\u00f76=57, 0\"],\n  [4, 2, \"432\u00f75=86, 2\", \"766\u00f77=109, 3\"],\n  [4, 3, \"880\u00f72=440, 0\", \"668\u00f73=222, 2\"],\n  [4, 4, \"622\u00f76=103, 4\", \"270\u00f72=135, 0\"],\n  [8, 0, \"482\u00f79=53, 5\", \"973\u00f77=139, 0\"],\n  [8, 1, \"984\u00f77=140, 4\", \"926\u00f78=115, 6\"],\n  [8, 2, \"543\u00f77=77, 4\", \"975\u00f74=243, 3\"],\n  [8, 3, \"176\u00f78=22, 0\", \"728\u00f75=145, 3\"],\n  [8, 4, \"908\u00f77=129, 5\", \"786\u00f75=157, 1\"],\n  [12, 0, \"546\u00f74=136, 2\", \"952\u00f79=105, 7\"],\n  [12, 1, \"517\u00f78=64, 5\", \"817\u00f75=163, 2\"],\n  [12, 2, \"782\u00f74=195, 2\", \"229\u00f76=38, 1\"],\n  [12, 3, \"527\u00f78=65, 7\", \"802\u00f76=133, 4\"],\n  [12, 4, \"776\u00f76=129, 2\", \"466\u00f73=155, 1\"],\n  [16, 0, \"289\u00f74=72, 1\", \"839\u00f72=419, 1\"],\n  [16, 1, \"567\u00f79=63, 0\", \"935\u00f74=233, 3\"],\n  [16, 2, \"204\u00f74=51, 0\", \"729\u00f78=91, 1\"],\n  [16, 3, \"338\u00f75=67, 3\", \"788\u00f74=197, 0\"],\n  [16, 4, \"494\u00f77=70, 4\", \"565\u00f72=282, 1\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Process each cell individually: search inside the cell body for its\n// exact current text and replace just that run's text in place, so run\n// formatting (font, size) and paragraph formatting (alignment) survive.\nfor (const [rowIndex, colIndex, oldText, newText] of replacements) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const searchResults = cell.body.search(oldText, { matchCase: true });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length > 0) {\n    searchResults.items[0].insertText(newText, Word.InsertLocation.replace);\n  } else {\n    // Fallback: if the expected old text isn't found (already changed,\n    // or formatting differs), just overwrite the cell body text plainly.\n    cell.body.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"divisor problem\" answers in the document's table with\n# their newly generated values, cell by cell, in document order.\n# Using direct Cell.Range.Text assignment (rather than Find/Replace) keeps\n# each edit strictly scoped to its own cell -- important because some of\n# the new values collide with old values used elsewhere in the table\n# (e.g. \"176\u00f78=22, 0\" is both an old value in one cell and a new value in\n# another), which a document-wide Find/Replace would handle incorrectly.\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$replacements = @(\n    @{ Row = 1;  Col = 1; Old = \"809\u00f79=89, 8\";   New = \"749\u00f76=124, 5\" },\n    @{ Row = 1;  Col = 2; Old = \"271\u00f77=38, 5\";   New = \"115\u00f79=12, 7\" },\n    @{ Row = 1;  Col = 3; Old = \"841\u00f78=105, 1\";  New = \"673\u00f74=168, 1\" },\n    @{ Row = 1;  Col = 4; Old = \"303\u00f76=50, 3\";   New = \"847\u00f79=94, 1\" },\n    @{ Row = 1;  Col = 5; Old = \"886\u00f74=221, 2\";  New = \"702\u00f76=117, 0\" },\n    @{ Row = 5;  Col = 1; Old = \"845\u00f79=93, 8\";   New = \"176\u00f78=22, 0\" },\n    @{ Row = 5;  Col = 2; Old = \"324\u00f73=108, 0\";  New = \"342\u00f76=57, 0\" },\n    @{ Row = 5;  Col = 3; Old = \"432\u00f75=86, 2\";   New = \"766\u00f77=109, 3\" },\n    @{ Row = 5;  Col = 4; Old = \"880\u00f72=440, 0\";  New = \"668\u00f73=222, 2\" },\n    @{ Row = 5;  Col = 5; Old = \"622\u00f76=103, 4\";  New = \"270\u00f72=135, 0\" },\n    @{ Row = 9;  Col = 1; Old = \"482\u00f79=53, 5\";   New = \"973\u00f77=139, 0\" },\n    @{ Row = 9;  Col = 2; Old = \"984\u00f77=140, 4\";  New = \"926\u00f78=115, 6\" },\n    @{ Row = 9;  Col = 3; Old = \"543\u00f77=77, 4\";   New = \"975\u00f74=243, 3\" },\n    @{ Row = 9;  Col = 4; Old = \"176\u00f78=22, 0\";   New = \"728\u00f75=145, 3\" },\n    @{ Row = 9;  Col = 5; Old = \"908\u00f77=129, 5\";  New = \"786\u00f75=157, 1\" },\n    @{ Row = 13; Col = 1; Old = \"546\u00f74=136, 2\";  New = \"952\u00f79=105, 7\" },\n    @{ Row = 13; Col = 2; Old = \"517\u00f78=64, 5\";   New = \"817\u00f75=163, 2\" },\n    @{ Row = 13; Col = 3; Old = \"782\u00f74=195, 2\";  New = \"229\u00f76=38, 1\" },\n    @{ Row = 13; Col = 4; Old = \"527\u00f78=65, 7\";   New = \"802\u00f76=133, 4\" },\n    @{ Row = 13; Col = 5; Old = \"776\u00f76=129, 2\";  New = \"466\u00f73=155, 1\" },\n    @{ Row = 17; Col = 1; Old = \"289\u00f74=72, 1\";   New = \"839\u00f72=419, 1\" },\n    @{ Row = 17; Col = 2; Old = \"567\u00f79=63, 0\";   New = \"935\u00f74=233, 3\" },\n    @{ Row = 17; Col = 3; Old = \"204\u00f74=51, 0\";   New = \"729\u00f78=91, 1\" },\n    @{ Row = 17; Col = 4; Old = \"338\u00f75=67, 3\";   New = \"788\u00f74=197, 0\" },\n    @{ Row = 17; Col = 5; Old = \"494\u00f77=70, 4\";   New = \"565\u00f72=282, 1\" }\n)\n\nforeach ($item in $replacements) {\n    $cell = $t.Cell($item.Row, $item.Col)\n    $rng = $cell.Range\n    # Assigning straight back onto the cell's Range.Text preserves the\n    # run/paragraph formatting already on that cell (unlike Find/Replace,\n    # which in this host searches/replaces document-wide rather than\n    # staying scoped to $rng, and unlike Range.Text = \"\" + insert, which\n    # would drop the rPr/pPr).\n    $rng.Text = $item.New\n}\n\nWrite-Output \"done\"\n"}
